$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain decimal numbers (e.g. "0.7393", "26.40").
# Force those cells to Text format before assigning so Excel keeps the exact
# string (matching "30.090" style values that already survive as text because
# they contain multiple dots and are never number-like).
$ws.Range("D2").Value = '30.090.49'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.909.53'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7393'
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.55'
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3089'
$ws.Range("E8").Value = '  -2.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.40'
$ws.Range("E9").Value = '  -5.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06966'
$ws.Range("E10").Value = '  -0.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08076'
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7695'
$ws.Range("E12").Value = '  -1.41%  '
$ws.Range("D13").Value = '1.903.25'
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.327'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.17'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.30'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '30.096.85'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.085'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007834'
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.91'
$ws.Range("E20").Value = '  -5.13%  '
$ws.Range("D21").Value = '2.190.21'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.078'
$ws.Range("E24").Value = '  +5.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.394'
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.20'
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.95'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1276'
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.054'
$ws.Range("E29").Value = '  -7.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.552'
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.336'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.077'
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.304'
$ws.Range("E34").Value = '  -0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05138'
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7479'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01955'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.799'
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.335'
$ws.Range("E40").Value = '  -4.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4502'
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.42'
$ws.Range("E42").Value = '  -5.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.981'
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8399'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.731'
$ws.Range("E46").Value = '  +1.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.964'
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.87'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = '2.079.87'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.75'
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1181'
$ws.Range("E51").Value = '  -3.90%  '
